$p = $ppt.ActivePresentation

# --- 1) Swap the table style on the three tables (slides 14, 15, 16) from the
#     deck's custom "Table_0" style to the built-in "No Style, Table Grid"
#     style, matching the Table Design gallery action in PowerPoint. ---
$newTableStyle = "{3BBA48EE-273E-43DC-8394-27CE8E00801F}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Re-colour the presentation's theme (Design > Themes) from the
#     custom "Integral" / "Red Violet" palette to the built-in "Office"
#     colour scheme. The format/font schemes are identical between the two
#     themes, so only the 12 theme colours need updating. ---
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
